$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 254.978138622258
$ws.Range("F3").Value = 359.692968632402
$ws.Range("E4").Value = 301.607225189755
$ws.Range("F4").Value = 448.894554900978
$ws.Range("E5").Value = 601.324735147575
$ws.Range("F5").Value = 1247.45976570351
$ws.Range("E6").Value = 768.928802655745
$ws.Range("F6").Value = 1438.17264810886
$ws.Range("B7").Value = 2.27381621367978
$ws.Range("E7").Value = 544.578837438132
$ws.Range("F7").Value = 1203.66923960347
$ws.Range("I7").Value = -3.72618378632022
$ws.Range("B8").Value = 128.666225377227
$ws.Range("E8").Value = 654.396529078582
$ws.Range("F8").Value = 1465.25206642596
$ws.Range("I8").Value = 59.6662253772271
$ws.Range("B9").Value = 340.465096986078
$ws.Range("E9").Value = 939.264048144212
$ws.Range("F9").Value = 1602.8287961124
$ws.Range("I9").Value = 271.465096986078
$ws.Range("B10").Value = 208.450482066848
$ws.Range("E10").Value = 796.223578999396
$ws.Range("F10").Value = 1487.26456139706
$ws.Range("I10").Value = 71.4504820668477
$ws.Range("B11").Value = 54.9824535964919
$ws.Range("E11").Value = 578.154910677943
$ws.Range("F11").Value = 1351.33633918474
$ws.Range("I11").Value = 12.9824535964919
$ws.Range("E12").Value = 520.782994662582
$ws.Range("F12").Value = 1132.60021207638
$ws.Range("F13").Value = 1043.06483211359
$ws.Range("E14").Value = 554.474684315649
$ws.Range("F14").Value = 1173.44828744567
$ws.Range("B15").Value = 0.0736080827840944
$ws.Range("E15").Value = 873.024726509974
$ws.Range("F15").Value = 1494.91606565094
$ws.Range("I15").Value = 0.0736080827840944
$ws.Range("B16").Value = 10.7961920147094
$ws.Range("E16").Value = 981.321326573523
$ws.Range("F16").Value = 1532.07993251643
$ws.Range("I16").Value = 10.7961920147094
$ws.Range("E17").Value = 965.346477532478
$ws.Range("F17").Value = 1610.29335679672
$ws.Range("E18").Value = 987.826132585542
$ws.Range("F18").Value = 1983.76800899424
$ws.Range("E19").Value = 914.170394632834
$ws.Range("F19").Value = 1694.63913104194
$ws.Range("B20").Value = 110.723154951442
$ws.Range("E20").Value = 851.723696579942
$ws.Range("F20").Value = 1663.0893061776
$ws.Range("I20").Value = 93.7231549514422
$ws.Range("B21").Value = 335.476116248648
$ws.Range("E21").Value = 1187.55630152838
$ws.Range("F21").Value = 1842.32015949287
$ws.Range("I21").Value = 263.476116248648
$ws.Range("B22").Value = 234.510692885332
$ws.Range("E22").Value = 1212.64994068625
$ws.Range("F22").Value = 1968.55517119565
$ws.Range("I22").Value = 171.510692885332
$ws.Range("B23").Value = 83.1130949616459
$ws.Range("E23").Value = 1195.46628780686
$ws.Range("F23").Value = 1921.90537147396
$ws.Range("I23").Value = 44.1130949616459
$ws.Range("B24").Value = 7.82597877361967
$ws.Range("E24").Value = 811.776356939282
$ws.Range("F24").Value = 1452.05990202354
$ws.Range("I24").Value = -6.17402122638033
$ws.Range("E25").Value = 757.313050655695
$ws.Range("F25").Value = 1625.06819155231
$ws.Range("E26").Value = 760.43129878723
$ws.Range("F26").Value = 1709.77872326811
$ws.Range("B27").Value = 0
$ws.Range("E27").Value = 1160.32257148498
$ws.Range("F27").Value = 1835.02068496425
$ws.Range("I27").Value = -1
$ws.Range("J27").Value = "Increase"
$ws.Range("B28").Value = 35.5985814042469
$ws.Range("E28").Value = 1410.37736156058
$ws.Range("F28").Value = 1935.46815587777
$ws.Range("I28").Value = 33.5985814042469
$ws.Range("B29").Value = 37.9130259817929
$ws.Range("E29").Value = 1120.84299270834
$ws.Range("F29").Value = 1962.92995275388
$ws.Range("I29").Value = 37.9130259817929
$ws.Range("E30").Value = 1208.16431370196
$ws.Range("F30").Value = 2036.62938102858
$ws.Range("E31").Value = 1042.38061536435
$ws.Range("F31").Value = 2091.3429109561
$ws.Range("B32").Value = 92.8317112176975
$ws.Range("E32").Value = 1082.73464097861
$ws.Range("F32").Value = 1862.29010932523
$ws.Range("I32").Value = 79.8317112176975
$ws.Range("B33").Value = 310.790561157857
$ws.Range("E33").Value = 1366.98797933007
$ws.Range("F33").Value = 2171.93339014986
$ws.Range("I33").Value = 253.790561157857
$ws.Range("B34").Value = 244.665092082335
$ws.Range("E34").Value = 1422.93610939132
$ws.Range("F34").Value = 2087.75855554226
$ws.Range("I34").Value = 191.665092082335
$ws.Range("B35").Value = 130.214616337173
$ws.Range("E35").Value = 1352.66907128059
$ws.Range("F35").Value = 2309.75513314184
$ws.Range("I35").Value = 103.214616337173
$ws.Range("B36").Value = 42.4495800982366
$ws.Range("E36").Value = 1045.41603705987
$ws.Range("F36").Value = 1929.55630275681
$ws.Range("I36").Value = 41.4495800982366
$ws.Range("E37").Value = 924.129813685913
$ws.Range("F37").Value = 1939.87490747641
$ws.Range("E38").Value = 900.144436576483
$ws.Range("F38").Value = 1820.43014040812
$ws.Range("E39").Value = 1386.67804516771
$ws.Range("F39").Value = 2185.96921927416
$ws.Range("B40").Value = 9.45860264206813
$ws.Range("E40").Value = 1496.43686123418
$ws.Range("F40").Value = 2202.05986232484
$ws.Range("I40").Value = 6.45860264206813
$ws.Range("B41").Value = 72.9849152134362
$ws.Range("E41").Value = 1486.034179702
$ws.Range("F41").Value = 2246.66894425637
$ws.Range("I41").Value = 71.9849152134362
$ws.Range("B42").Value = 35.8347305774807
$ws.Range("E42").Value = 1293.06371296081
$ws.Range("F42").Value = 2097.27962185932
$ws.Range("I42").Value = 35.8347305774807
$ws.Range("B43").Value = 1.05504620735705
$ws.Range("E43").Value = 1360.47991576257
$ws.Range("F43").Value = 2316.48569727922
$ws.Range("I43").Value = -1.94495379264295
$ws.Range("B44").Value = 90.2476490828214
$ws.Range("E44").Value = 1184.80528973202
$ws.Range("F44").Value = 2011.68402622186
$ws.Range("I44").Value = 57.2476490828214
$ws.Range("B45").Value = 286.771438411819
$ws.Range("E45").Value = 1514.99579073141
$ws.Range("F45").Value = 2099.55985156725
$ws.Range("I45").Value = 206.771438411819
$ws.Range("B46").Value = 200.503894375854
$ws.Range("E46").Value = 1560.04261073404
$ws.Range("F46").Value = 2191.24604102313
$ws.Range("I46").Value = 138.503894375854
$ws.Range("B47").Value = 153.999406877275
$ws.Range("E47").Value = 1443.50524495038
$ws.Range("F47").Value = 2311.75722156458
$ws.Range("I47").Value = 121.999406877275
$ws.Range("B48").Value = 78.0042928001795
$ws.Range("E48").Value = 1356.38925532798
$ws.Range("F48").Value = 2292.55225304514
$ws.Range("I48").Value = 66.0042928001795
$ws.Range("B49").Value = 11.5414797704154
$ws.Range("E49").Value = 1260.28741822718
$ws.Range("F49").Value = 2071.95209271933
$ws.Range("I49").Value = 7.54147977041536
$ws.Range("J49").Value = "Decrease"
